$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new observation rows (13-15) appended below the existing data.
$rows = @(
    @{ Row=13; A=131063376; Q=397688; R=7048558; Z="16:02"; AB="16:02" },
    @{ Row=14; A=131061121; Q=397980; R=7048389; Z="13:22"; AB="13:22" },
    @{ Row=15; A=131060978; Q=397981; R=7048398; Z="13:13"; AB="13:13" }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Cells.Item($n, 1).Value = $r.A          # A  Id
    $ws.Cells.Item($n, 2).Value = 91828          # B  Taxonsorteringsordning
    $ws.Cells.Item($n, 4).Value = "NT"           # D  Rodlistade
    $ws.Cells.Item($n, 5).Value = 5432           # E  TaxonId
    $ws.Cells.Item($n, 6).Value = "Granticka"    # F  Artnamn
    $ws.Cells.Item($n, 7).Value = "Porodaedalea chrysoloma s.lat."  # G Vetenskapligt namn

    # H, I are present in the source export as empty (text) placeholder
    # cells; Excel has no way to persist an empty-but-typed cell other
    # than leaving it unset, which yields the same effective blank value.
    $ws.Cells.Item($n, 8).Value = ""
    $ws.Cells.Item($n, 9).Value = ""

    $ws.Cells.Item($n, 16).Value = "Vallen, Vallen, Jmt"  # P Lokalnamn
    $ws.Cells.Item($n, 17).Value = $r.Q                    # Q Ost
    $ws.Cells.Item($n, 18).Value = $r.R                    # R Nord
    $ws.Cells.Item($n, 19).Value = 10                      # S Noggrannhet
    $ws.Cells.Item($n, 20).Value = "Jämtland"               # T Lan
    $ws.Cells.Item($n, 21).Value = "Åre"                    # U Kommun
    $ws.Cells.Item($n, 22).Value = "Jämtland"               # V Provins
    $ws.Cells.Item($n, 23).Value = "Kall"                   # W Socken

    # Y, Z, AA, AB are plain-text date/time strings in the source data.
    # Pre-formatting the cell as Text keeps Excel from coercing the
    # "yyyy-mm-dd" strings into date serial numbers.
    $ws.Cells.Item($n, 25).NumberFormat = "@"
    $ws.Cells.Item($n, 25).Value = "2026-02-07"   # Y Startdatum
    $ws.Cells.Item($n, 25).Style = "Normal"

    $ws.Cells.Item($n, 26).NumberFormat = "@"
    $ws.Cells.Item($n, 26).Value = $r.Z           # Z Starttid
    $ws.Cells.Item($n, 26).Style = "Normal"

    $ws.Cells.Item($n, 27).NumberFormat = "@"
    $ws.Cells.Item($n, 27).Value = "2026-02-07"   # AA Slutdatum
    $ws.Cells.Item($n, 27).Style = "Normal"

    $ws.Cells.Item($n, 28).NumberFormat = "@"
    $ws.Cells.Item($n, 28).Value = $r.AB          # AB Sluttid
    $ws.Cells.Item($n, 28).Style = "Normal"

    $ws.Cells.Item($n, 30).Value = $false   # AD Ej aterfunnen
    $ws.Cells.Item($n, 31).Value = $false   # AE Osaker artbestamning
    $ws.Cells.Item($n, 33).Value = $false   # AG Ospontan

    $ws.Cells.Item($n, 46).Value = ""   # AT Bestamningsar (empty)

    $ws.Cells.Item($n, 49).Value = "Fabian Pettersson"  # AW Rapportor
    $ws.Cells.Item($n, 50).Value = "Fabian Pettersson"  # AX Observatorer

    $ws.Cells.Item($n, 51).Value = ""   # AY Projektnamn (empty)
}
